$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.802.08"
$ws.Range("E2").Value = "  -5.40%  "

# Row 3
$ws.Range("D3").Value = "2.209.77"
$ws.Range("E3").Value = "  -7.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "98.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -9.71%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.582"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -7.96%  "

# Row 8
$ws.Range("E8").Value = "  +0.07%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.558"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -9.88%  "

# Row 10
$ws.Range("E10").Value = "  -11.56%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.20%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0824"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -10.37%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -9.04%  "

# Row 14
$ws.Range("E14").Value = "  -3.94%  "

# Row 15
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.547.38"
$ws.Range("E15").Value = "  -6.98%  "

# Row 16
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.858"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -12.55%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.82%  "

# Row 18
$ws.Range("D18").Value = "2.212.09"
$ws.Range("E18").Value = "  -6.60%  "

# Row 19
$ws.Range("D19").Value = "42.698.21"
$ws.Range("E19").Value = "  -5.58%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.79"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.48%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0957"
$ws.Range("E21").Value = "  -9.87%  "

# Row 22
$ws.Range("E22").Value = "  -12.84%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -11.40%  "

# Row 24
$ws.Range("E24").Value = "  -10.00%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "235.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.64%  "

# Row 26
$ws.Range("E26").Value = "  -8.47%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.03%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.62%  "

# Row 29
$ws.Range("E29").Value = "  -5.35%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.26"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -13.94%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.08%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0868"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -10.36%  "

# Row 33
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "33.76"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -10.24%  "

# Row 34
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.42"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.54%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.70%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.17"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.56%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.91%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.121"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.82%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.43%  "

# Row 40
$ws.Range("E40").Value = "  -12.71%  "

# Row 41
$ws.Range("E41").Value = "  -7.05%  "

# Row 42
$ws.Range("E42").Value = "  -8.69%  "

# Row 43
$ws.Range("D43").Value = "1.881.27"
$ws.Range("E43").Value = "  +0.15%  "

# Row 44
$ws.Range("E44").Value = "  +0.28%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.13"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.49%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -11.39%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.205"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -10.46%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.37"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.05%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "60.22"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -13.61%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.59"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -8.60%  "

# Row 51
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -7.58%  "
